$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Caja"

# Update "Donaciones" total (row 4)
$ws.Range("B4").Value = 4476000

# Update "Total Ingresos" (row 7)
$ws.Range("B7").Value = 4488500

# Update "aaa" (row 9)
$ws.Range("B9").Value = 13428

# Row 10 used to be "xd" with values 1000/0/0/0; it becomes "Total Egresos"
# with values 13428/0/0/0, styled like the other total rows (gray fill + border).
$ws.Range("A10").Value = "Total Egresos"
$ws.Range("B10").Value = 13428
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("B7:E7").Copy()
$ws.Range("B10:E10").PasteSpecial(-4122)

# Row 11 used to be "Hola" with values 11200/0/0/0; it becomes "Acumulado"
# with values 4475072/0/0/0, styled like the other total rows (gray fill + border).
$ws.Range("A11").Value = "Acumulado"
$ws.Range("B11").Value = 4475072
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("B7:E7").Copy()
$ws.Range("B11:E11").PasteSpecial(-4122)

# Delete the old rows 12 and 13 ("Total Estimaciones" and "Acumulado")
$ws.Range("A12:E13").EntireRow.Delete()

# Narrow column A (ColumnWidth uses character units; the stored XML "width"
# attribute is ColumnWidth + ~0.83 due to Excel's internal padding formula,
# so set 15.17 to land on a serialized width of 16).
$ws.Columns.Item(1).ColumnWidth = 15.17
